# CTEMMM.xlsx maintenance update
#   - Print_Area adjustments on both sheets (room for a footer/copyright row)
#   - New copyright / license footer row on 'CTEM-MM'
#   - Page setup tweaks (print scale, Definitions sheet fit-to-page + orientation)
#   - Definitions sheet selection anchor follows its new print area

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CTEM-MM")
$ws2 = $wb.Worksheets.Item("Definitions")

# --- 'CTEM-MM': add the footer/copyright notice below the table ------------
$ws1.Range("B44").Value = "© 2025 ZenzizenSec Inc. All rights reserved. CTEMMM™ – Continuous Threat Exposure Management Maturity Model. Licensed for non-commercial use only. For terms, contact license@zenzizensec.com"
$ws1.Range("B44").Font.Name = "Times New Roman"
$ws1.Range("B44").Font.Size = 14
$ws1.Range("B44").HorizontalAlignment = -4131
$ws1.Range("B44").VerticalAlignment = -4108
$ws1.Range("B44").WrapText = $false

# --- Print areas: extend 'CTEM-MM' to cover the new footer row, and define one for 'Definitions'
$ws1.PageSetup.PrintArea = '$A$3:$H$44'
$ws2.PageSetup.PrintArea = '$A$2:$C$45'

# --- 'CTEM-MM' print scale tweak (40% -> 39%) so the wider area still fits a page
$ws1.PageSetup.Zoom = 39
$ws1.PageSetup.FitToPagesTall = 5

# --- 'Definitions' gains fit-to-page printing (portrait, 2 pages tall, 65% scale)
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.Zoom = 65
$ws2.PageSetup.FitToPagesWide = 1
$ws2.PageSetup.FitToPagesTall = 2

# --- keep the 'Definitions' selection in sync with its new print area
$ws2.Range("A2:C45").Select() | Out-Null
